$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'25.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.021"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05621"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.570"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Value = "'0.8140"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8366"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1338"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Value = "'0.02841"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.09400"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.001507"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.0005944"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'14OneONE"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006110"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'3.500"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Value = "'0.03264"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.1291"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'3.750"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.04686"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Value = "'0.004532"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.00009695"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'26NitroExNTXBestin24h"
$ws.Range("E27").Style = "Normal"
$ws.Range("B41").Value = "'BKEXToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.1361"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'40BKEXTokenBKK"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'CEJI"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.002733"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'41CEJICEJI"
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'KickToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.003382"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'42KickTokenKICKWorstin24h"
$ws.Range("E43").Style = "Normal"
$ws.Range("D45").Value = "'0.00005289"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.2258"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'46CoinbaseStockTokenCOIN"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002020"
$ws.Range("D48").Style = "Normal"
